$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New block appended for the 15-10-2020 snapshot (rows 1009-1044) ---

# Row 1009 is the periodic "States/UT" column-header row that repeats every 36 rows
# throughout this sheet (see rows 1, 37, 73, ... 973). Reproduce it, with its header
# formatting (bold, thin-bordered, centered), by copying the most recent such header row.
$ws.Range("A973:H973").Copy($ws.Range("A1009:H1009"))

# Rows 1010-1044: one data row per State/UT for 15-10-2020, in the same column order as
# the rest of the sheet: States/UT, Active Cases, Active Cases Since Yesterday,
# Recovered Cases, Recovered Cases Since Yesterday, Deceased Cases,
# Deceased Cases Since Yesterday, Date.
$newData = @(
    @('Andaman and Nicobar Islands', 195, -4, 3796, 14, 55, 0, '15-10-2020'),
    @('Andhra Pradesh', 41669, -1186, 719477, 5050, 6319, 28, '15-10-2020'),
    @('Arunachal Pradesh', 3045, 85, 9694, 121, 29, 1, '15-10-2020'),
    @('Assam', 29307, 410, 168072, 1013, 834, 4, '15-10-2020'),
    @('Bihar', 10756, -79, 188380, 1321, 967, 6, '15-10-2020'),
    @('Chandigarh', 1085, -42, 12119, 112, 199, 2, '15-10-2020'),
    @('Chhattisgarh', 27809, 601, 121548, 2196, 1339, 33, '15-10-2020'),
    @('Dadra and Nagar Haveli and Daman and Diu', 79, -6, 3090, 9, 2, 0, '15-10-2020'),
    @('Delhi', 21903, 413, 289747, 2867, 5898, 44, '15-10-2020'),
    @('Goa', 4188, -128, 34731, 479, 519, 5, '15-10-2020'),
    @('Gujarat', 14937, -250, 136404, 1414, 3595, 11, '15-10-2020'),
    @('Haryana', 10187, -132, 133706, 1324, 1614, 13, '15-10-2020'),
    @('Himachal Pradesh', 2520, 13, 15233, 232, 255, 1, '15-10-2020'),
    @('Jammu and Kashmir', 9739, -127, 74318, 816, 1352, 12, '15-10-2020'),
    @('Jharkhand', 7191, -426, 86367, 1053, 811, 6, '15-10-2020'),
    @('Karnataka', 114006, 528, 611167, 8662, 10198, 75, '15-10-2020'),
    @('Kerala', 93925, -1568, 215149, 7792, 1066, 20, '15-10-2020'),
    @('Ladakh', 979, 10, 4261, 56, 64, 0, '15-10-2020'),
    @('Madhya Pradesh', 14432, -229, 138158, 5729, 2686, 15, '15-10-2020'),
    @('Maharashtra', 196761, -9123, 1316769, 19517, 40859, 158, '15-10-2020'),
    @('Manipur', 3076, 209, 10915, 86, 103, 6, '15-10-2020'),
    @('Meghalaya', 2339, -28, 5582, 176, 70, 5, '15-10-2020'),
    @('Mizoram', 112, -7, 2108, 15, 0, 0, '15-10-2020'),
    @('Nagaland', 1478, -35, 5916, 85, 22, 3, '15-10-2020'),
    @('Odisha', 22716, -176, 235763, 2775, 1062, 5, '15-10-2020'),
    @('Puducherry', 4525, -47, 27152, 287, 568, 1, '15-10-2020'),
    @('Punjab', 7760, -452, 114075, 970, 3925, 31, '15-10-2020'),
    @('Rajasthan', 21711, -213, 141835, 2219, 1694, 15, '15-10-2020'),
    @('Sikkim', 325, -19, 3075, 61, 59, 0, '15-10-2020'),
    @('Tamil Nadu', 42566, -673, 617403, 5083, 10423, 52, '15-10-2020'),
    @('Telengana', 23203, -525, 193218, 1949, 1249, 8, '15-10-2020'),
    @('Tripura', 3318, -182, 25394, 353, 319, 1, '15-10-2020'),
    @('Uttarakhand', 6145, -431, 49129, 846, 796, 14, '15-10-2020'),
    @('Uttar Pradesh', 36898, -1184, 401306, 3736, 6507, 41, '15-10-2020'),
    @('West Bengal', 31505, 517, 268384, 3096, 5808, 64, '15-10-2020')
)

$r = 1010
foreach ($row in $newData) {
    $ws.Range("A$r").Value() = $row[0]
    $ws.Range("B$r").Value() = $row[1]
    $ws.Range("C$r").Value() = $row[2]
    $ws.Range("D$r").Value() = $row[3]
    $ws.Range("E$r").Value() = $row[4]
    $ws.Range("F$r").Value() = $row[5]
    $ws.Range("G$r").Value() = $row[6]
    $ws.Range("H$r").Value() = $row[7]
    $r = $r + 1
}
